$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 "HSI_wet_AICc" ---
# New "delta" header column
$ws1.Range("D4").Value = "delta"

# Number format for the AICc column (column C)
$ws1.Range("C5:C11").NumberFormat = "0.00"

# delta = AICc - 427.58 (best model's AICc), replacing the old "* best" label in D9
$ws1.Range("D5").Formula = "=C5-427.58"
$ws1.Range("D6:D11").Formula = "=C6-427.58"
$ws1.Range("D10").ClearContents()
$ws1.Range("D5:D11").NumberFormat = "0.00"

# --- Sheet2 "liverFA R results" ---
$ws2.Range("I12").Value = "too much skew in data"
$ws2.Range("I13").Value = "over fit. Did not use model 16"

# --- View state: sheet1 becomes the active/selected tab ---
$ws2.Range("H14").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Activate()
$ws1.Range("D10").Select()
